$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the phone numbers / office-hours text for the existing Elementary
# teachers (columns D and E of rows 1-4).
# ---------------------------------------------------------------------------

$ws.Range("D1").Value = "(239) 281-4712"
$ws.Range("E1").Value = "Thurs 9am-4pm"

$ws.Range("D2").Value = "(239) 822-0318"
$ws.Range("E2").Value = "Wed 1pm-3pm, Thur 10am-12pm"

$ws.Range("D3").Value = "(239) 357-7602"
$ws.Range("E3").Value = "Please make an appointment"

$ws.Range("D4").Value = "(239) 823-8065"
$ws.Range("E4").Value = "Mon- Frid 8am-8pm"

# ---------------------------------------------------------------------------
# Add the new "Secondary" teacher rows (5-8).
# ---------------------------------------------------------------------------

$ws.Range("A5").Value = "Secondary"
$ws.Range("B5").Value = "Mr. Bruni"
$ws.Range("C5").Value = "Jamesdb@leeschools.net"
$ws.Range("D5").Value = "(239) 357-8038"
$ws.Range("E5").Value = "Mon 9am-12pm"

$ws.Range("A6").Value = "Secondary"
$ws.Range("B6").Value = "Mr. McKinley"
$ws.Range("C6").Value = "edwardamck@leeschools.net"
$ws.Range("D6").Value = "(239) 822-9725"
$ws.Range("E6").Value = "Thur 3pm-4pm"

$ws.Range("A7").Value = "Secondary"
$ws.Range("B7").Value = "Mrs. Ausman"
$ws.Range("C7").Value = "MeaganMA@leeschools.net"
$ws.Range("D7").Value = "(239) 823-2565"
$ws.Range("E7").Value = "Mon-Fri 7am-8pm"

$ws.Range("A8").Value = "Secondary"
$ws.Range("B8").Value = "Mrs. Brooks"
$ws.Range("C8").Value = "brittanypb@leeschools.net"
$ws.Range("D8").Value = "(239) 357-2709"
$ws.Range("E8").Value = "Thur 1pm-4pm by Appointment"

# ---------------------------------------------------------------------------
# Formatting. Two custom looks are used throughout the table:
#   - plain black Calibri 11 (font color forced to black / Automatic)
#   - Arial 10 (smaller font, used for the phone-number column and for all
#     of the new "Secondary" rows except the last phone number cell)
#
# Build each look once on a seed cell, then use Copy / PasteSpecial (paste
# formats only) to fan it out to the rest of the affected cells - this
# mirrors applying the look via the Format Painter and keeps the workbook
# from accumulating a separate style per cell.
# ---------------------------------------------------------------------------

$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Color = 0
$ws.Range("B1").Copy()
$ws.Range("A1:C4").PasteSpecial(-4122)
$ws.Range("E1:E4").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)

$ws.Range("D1").Font.Name = "Arial"
$ws.Range("D1").Font.Size = 10
$ws.Range("D1").Copy()
$ws.Range("D2:D4").PasteSpecial(-4122)
$ws.Range("A5:E7").PasteSpecial(-4122)
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Range("E8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("E13").Select()
